$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update "Sheet3" (UsersData's sibling tab) row1 A1 text value:
#    "Test launch" -> "Test launch1234"
# ------------------------------------------------------------------
$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3.Range("A1").Value = "Test launch1234"

# ------------------------------------------------------------------
# 2. Add a brand new worksheet "Sheet4" as the last tab of the workbook
#    and populate row 1 with the new test-case data.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet4 = $wb.Worksheets.Add($null, $lastSheet)
$sheet4.Name = "Sheet4"

$sheet4.Range("A1").Value = "selenium"
$sheet4.Range("B1").Value = "selenium@gmail.com"
$sheet4.Range("C1").Value = "apartment"
$sheet4.Range("D1").Value = "looking for apartment"
$sheet4.Range("E1").Value = 40000
$sheet4.Range("F1").Value = 2000
$sheet4.Range("G1").Value = 2
$sheet4.Range("H1").Value = 5

# Mailto hyperlink on B1 (matches the pattern used on the other sheets).
$sheet4.Hyperlinks.Add($sheet4.Range("B1"), "mailto:selenium@gmail.com")
$sheet4.Range("B1").Style = "Hyperlink"

# Final selection left on the new sheet is H1.
$sheet4.Range("H1").Select()

# ------------------------------------------------------------------
# 3. Make "UsersData" the active/selected tab of the workbook again.
# ------------------------------------------------------------------
$usersData = $wb.Worksheets.Item("UsersData")
$usersData.Activate()
